$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# Delete rows 3 and 4 (old SamplesTab/FilesTab rows), shifting rows up
$ws.Range("A3:G4").Delete()

# Update D2 and E2 filenames for TC03 test
$ws.Range("D2").Value = "TC03_Canine_E2E_MultipleFilters-Study_Breed_Sex_Neo4jData.xlsx"
$ws.Range("E2").Value = "TC03_Canine_E2E_MultipleFilters-Study_Breed_Sex_WebData.xlsx"

# Update selection/view
$ws.Range("C7").Select() | Out-Null

Write-Output "done"
